$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.Style = 'Normal'
}

# Row 2
Set-TextCell 2 4 '44.400.92'
Set-TextCell 2 5 '  +0.46%  '

# Row 3
Set-TextCell 3 4 '2.243.62'
Set-TextCell 3 5 '  -0.47%  '

# Row 4
Set-TextCell 4 5 '  +0.39%  '

# Row 5
Set-TextCell 5 4 '305.56'
Set-TextCell 5 5 '  -0.64%  '

# Row 6
Set-TextCell 6 4 '92.92'

# Row 7
Set-TextCell 7 5 '  -1.01%  '

# Row 8
Set-TextCell 8 5 '  +0.28%  '

# Row 9
Set-TextCell 9 4 '0.521'
Set-TextCell 9 5 '  -2.81%  '

# Row 10
Set-TextCell 10 4 '34.46'
Set-TextCell 10 5 '  -3.15%  '

# Row 11
Set-TextCell 11 4 '0.0809'
Set-TextCell 11 5 '  -1.72%  '

# Row 12
Set-TextCell 12 5 '  -2.97%  '

# Row 13
Set-TextCell 13 5 '  -0.34%  '

# Row 14
Set-TextCell 14 4 '2.364.56'
Set-TextCell 14 5 '  +3.27%  '

# Row 15
Set-TextCell 15 4 '0.835'
Set-TextCell 15 5 '  -0.72%  '

# Row 16
Set-TextCell 16 5 '  -2.27%  '

# Row 17
Set-TextCell 17 4 '44.077.14'
Set-TextCell 17 5 '  +0.06%  '

# Row 18
Set-TextCell 18 5 '  -1.91%  '

# Row 19
Set-TextCell 19 4 '12.30'
Set-TextCell 19 5 '  -4.57%  '

# Row 20
Set-TextCell 20 5 '  -0.13%  '

# Row 21
Set-TextCell 21 4 '65.54'
Set-TextCell 21 5 '  +0.02%  '

# Row 22
Set-TextCell 22 2 'BitcoinCash'
Set-TextCell 22 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 22 4 '237.09'
Set-TextCell 22 5 '  -2.36%  '

# Row 23
Set-TextCell 23 2 'PancakeSwap'
Set-TextCell 23 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 23 4 '2.93'
Set-TextCell 23 5 '  -0.76%  '

# Row 24
Set-TextCell 24 4 '1.97'
Set-TextCell 24 5 '  -0.23%  '

# Row 25
Set-TextCell 25 5 '  -0.01%  '

# Row 26
Set-TextCell 26 4 '38.80'
Set-TextCell 26 5 '  +5.20%  '

# Row 27
Set-TextCell 27 4 '2.20'
Set-TextCell 27 5 '  +2.76%  '

# Row 28
Set-TextCell 28 5 '  -3.43%  '

# Row 29
Set-TextCell 29 4 '5.92'
Set-TextCell 29 5 '  -4.20%  '

# Row 30
Set-TextCell 30 4 '20.01'
Set-TextCell 30 5 '  -0.68%  '

# Row 31
Set-TextCell 31 4 '153.80'
Set-TextCell 31 5 '  -2.23%  '

# Row 32
Set-TextCell 32 4 '0.0797'
Set-TextCell 32 5 '  -3.83%  '

# Row 33
Set-TextCell 33 5 '  -0.40%  '

# Row 34
Set-TextCell 34 5 '  -14.21%  '

# Row 35
Set-TextCell 35 5 '  +0.24%  '

# Row 36
Set-TextCell 36 5 '  -0.02%  '

# Row 37
Set-TextCell 37 5 '  -1.85%  '

# Row 38
Set-TextCell 38 4 '3.43'
Set-TextCell 38 5 '  +0.99%  '

# Row 39
Set-TextCell 39 4 '14.50'
Set-TextCell 39 5 '  -5.06%  '

# Row 40
Set-TextCell 40 4 '3.80'
Set-TextCell 40 5 '  -2.59%  '

# Row 41
Set-TextCell 41 4 '0.0300'
Set-TextCell 41 5 '  -2.10%  '

# Row 42
Set-TextCell 42 5 '  +0.43%  '

# Row 43
Set-TextCell 43 4 '1.729.13'
Set-TextCell 43 5 '  -1.81%  '

# Row 44
Set-TextCell 44 4 '0.193'
Set-TextCell 44 5 '  -0.23%  '

# Row 45
Set-TextCell 45 4 '80.35'
Set-TextCell 45 5 '  -9.36%  '

# Row 46
Set-TextCell 46 4 '99.17'
Set-TextCell 46 5 '  -2.45%  '

# Row 47
Set-TextCell 47 4 '4.92'
Set-TextCell 47 5 '  -4.83%  '

# Row 48
Set-TextCell 48 2 'MultiversX'
Set-TextCell 48 3 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextCell 48 4 '55.32'
Set-TextCell 48 5 '  -0.27%  '

# Row 49
Set-TextCell 49 4 '8.15'
Set-TextCell 49 5 '  -1.37%  '

# Row 50
Set-TextCell 50 2 'ordi'
Set-TextCell 50 3 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextCell 50 4 '69.54'
Set-TextCell 50 5 '  -1.20%  '

# Row 51
Set-TextCell 51 2 'Stacks'
Set-TextCell 51 3 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell 51 4 '1.58'
Set-TextCell 51 5 '  +1.80%  '
